# Update "Salario Basico" values (column G) for rows 16-18 on Hoja1
# from 781242 to 737717, reflecting the refreshed EC (estado de cuenta)
# database values referenced in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("G16").Value = 737717
$ws.Range("G17").Value = 737717
$ws.Range("G18").Value = 737717
